# Rename worksheets (new randomized task-order identifiers)
$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item(1).Name = "GNG_TO-16502912626227846"
$wb.Worksheets.Item(2).Name = "NB_TO-16502912657375357"
$wb.Worksheets.Item(3).Name = "RS_TO-1650291265738474"
$wb.Worksheets.Item(4).Name = "TOL_TO-16502912657857673"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16502912658786802"

# Sheet 1 (GNG) - update stim filenames in column B
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-1650291262568224.csv"
$ws1.Range("B3").Value = "GNG_stims-16502912625910232.csv"
$ws1.Range("B4").Value = "go_stims-16502912625930326.csv"
$ws1.Range("B5").Value = "GNG_stims-16502912626217852.csv"

# Sheet 2 (NB) - update stim filenames in column B
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "ZB-match_6-16502912634965765.csv"
$ws2.Range("B3").Value = "TB-16502912650122123.csv"
$ws2.Range("B4").Value = "OB-16502912647515936.csv"
$ws2.Range("B5").Value = "ZB-match_2-16502912632942007.csv"
$ws2.Range("B6").Value = "TB-1650291265115218.csv"
$ws2.Range("B7").Value = "ZB-match_0-1650291264204073.csv"
$ws2.Range("B8").Value = "TB-16502912657252526.csv"
$ws2.Range("B9").Value = "OB-16502912647918277.csv"
$ws2.Range("B10").Value = "OB-16502912646887205.csv"

# Sheet 4 (TOL) - update stim filenames in column B
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16502912657529967.csv"
$ws4.Range("B3").Value = "ZM_stims-16502912657404704.csv"
$ws4.Range("B4").Value = "MM_stims-1650291265768922.csv"
$ws4.Range("B5").Value = "ZM_stims-16502912657539976.csv"
$ws4.Range("B6").Value = "MM_stims-16502912657844229.csv"
$ws4.Range("B7").Value = "ZM_stims-1650291265768922.csv"

# Sheet 5 (vSAT) - update stim filenames in column B
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "vSAT_stims-16502912658630412.csv"
$ws5.Range("B3").Value = "vSAT_stims-16502912658312628.csv"
$ws5.Range("B4").Value = "SAT_stims-16502912657905576.csv"
$ws5.Range("B5").Value = "SAT_stims-16502912658150175.csv"
